$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the questionnaire script text in A16 with the new script text
$ws.Range("A16").Value = "Appeal information`nSite accessibility`nConservation areas`nMonuments and listed buildings`nPlans`nPlanning Officer's Report`nInterested Parties and Consultees`nInterested Parties and Consultee Comments`nLocal policy`nSuggested Conditions`nAdditional Information`nConfirmation"

# Adjust the row height to match the new (shorter) text block
$ws.Rows.Item(16).RowHeight = 174

# Update the view/selection state
$ws.Application.ActiveWindow.ScrollRow = 15
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("A16").Select()
